# "frame work data setup"
#
# Populate a second test-suite data row (TC_Name / Destination) on Sheet1,
# widen column A so the longer TC_Name text fits, and restore the window
# position/size that was in effect when the data was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of framework data.
$ws.Range("A2").Value = "enterDetailsHomePage"
$ws.Range("B2").Value = "Mum"

# Column A must widen to fit "enterDetailsHomePage". Range.ColumnWidth is in
# characters (Normal-style font) and Excel stores/rounds it on a whole-pixel
# grid, so feed it the character width whose rounded OOXML width lands as
# close as possible to the authored 32.42578125.
$ws.Columns.Item(1).ColumnWidth = (32.42578125 - 5/6)

# The cursor was left on B1 after the data entry.
$ws.Range("B1").Select()

# Restore the workbook window's last position/size.
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 2520
$win.Width = 11310
$win.Height = 3975
